$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 21.01.2022 18:30"

# Update D7 to a numeric value (was stored as text "+0.2")
$ws.Range("D7").Value = 0.2

# Update E7 to a numeric date/time value (was stored as text "2022-01-21 18:15:22").
# Set the number format first so the existing date-format style is reused
# instead of a brand new (unused) style being created.
$ws.Range("E7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E7").Value = (Get-Date -Year 2022 -Month 1 -Day 21 -Hour 18 -Minute 15 -Second 22)
